# Weekly update: a new price record for "Alcachofa" (Vega Modelo de Temuco)
# is inserted as row 94, pushing all subsequent records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 94 (shifts rows 94-161 down to 95-162,
# and extends the used range / dimension to A1:R162 automatically).
$ws.Rows.Item(94).EntireRow.Insert()

# Populate the newly inserted row 94 with the new record's data.
$ws.Cells.Item(94, 1).Value = 10
$ws.Cells.Item(94, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(94, 3).Value = "La Araucanía"
$ws.Cells.Item(94, 4).Value = 44719
$ws.Cells.Item(94, 5).Value = 9
$ws.Cells.Item(94, 6).Value = 100112013
$ws.Cells.Item(94, 7).Value = "Alcachofa"
$ws.Cells.Item(94, 8).Value = "Madrigal"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 40
$ws.Cells.Item(94, 11).Value = 23000
$ws.Cells.Item(94, 12).Value = 23000
$ws.Cells.Item(94, 13).Value = 23000
$ws.Cells.Item(94, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(94, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(94, 16).Value = 767
$ws.Cells.Item(94, 17).Value = 30
$ws.Cells.Item(94, 18).Value = "Hortaliza"
